$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "etststs"
$ws.Range("B3").Value = " "
